# Rewrites the LR-pairs data block (rows 2..10) with the updated TPM-derived values.
# Row 1 (headers) is untouched; columns A-T keep the same schema.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object "object[,]" 9,20
# row 2: ECs -> ECs
$data[0,0] = 'ECs'
$data[0,1] = 'Tnfsf10'
$data[0,2] = 'Tnfrsf11b'
$data[0,3] = 'ECs'
$data[0,4] = 3
$data[0,5] = 1
$data[0,6] = 13.86674266666667
$data[0,7] = 41.600228
$data[0,8] = 0.9894894203110381
$data[0,9] = 0.989489420311038
$data[0,10] = 1
$data[0,11] = 0.3333333333333333
$data[0,12] = 0.07580833333333332
$data[0,13] = 0.227425
$data[0,14] = 0.028190957994264
$data[0,15] = 0.02819095799426401
$data[0,16] = 1.051214650322222
$data[0,17] = 9.4609318529
$data[0,18] = 0.02789465468375711
$data[0,19] = 0.02789465468375711
# row 3: ECs -> FAPs
$data[1,0] = 'ECs'
$data[1,1] = 'Tnfsf10'
$data[1,2] = 'Tnfrsf11b'
$data[1,3] = 'FAPs'
$data[1,4] = 3
$data[1,5] = 1
$data[1,6] = 13.86674266666667
$data[1,7] = 41.600228
$data[1,8] = 0.9894894203110381
$data[1,9] = 0.989489420311038
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 2.095195666666667
$data[1,13] = 6.285587
$data[1,14] = 0.7791435378093522
$data[1,15] = 0.7791435378093522
$data[1,16] = 29.05353914598178
$data[1,17] = 261.481852313836
$data[1,18] = 0.7709542875660673
$data[1,19] = 0.7709542875660672
# row 4: ECs -> MuSCs
$data[2,0] = 'ECs'
$data[2,1] = 'Tnfsf10'
$data[2,2] = 'Tnfrsf11b'
$data[2,3] = 'MuSCs'
$data[2,4] = 3
$data[2,5] = 1
$data[2,6] = 13.86674266666667
$data[2,7] = 41.600228
$data[2,8] = 0.9894894203110381
$data[2,9] = 0.989489420311038
$data[2,10] = 3
$data[2,11] = 1
$data[2,12] = 0.5180969999999999
$data[2,13] = 1.554291
$data[2,14] = 0.1926655041963838
$data[2,15] = 0.1926655041963838
$data[2,16] = 7.184317775371999
$data[2,17] = 64.65885997834799
$data[2,18] = 0.1906404780612137
$data[2,19] = 0.1906404780612137
# row 5: FAPs -> ECs
$data[3,0] = 'FAPs'
$data[3,1] = 'Tnfsf10'
$data[3,2] = 'Tnfrsf11b'
$data[3,3] = 'ECs'
$data[3,4] = 1
$data[3,5] = 0.3333333333333333
$data[3,6] = 0.09171866666666667
$data[3,7] = 0.275156
$data[3,8] = 0.006544770642485517
$data[3,9] = 0.006544770642485517
$data[3,10] = 1
$data[3,11] = 0.3333333333333333
$data[3,12] = 0.07580833333333332
$data[3,13] = 0.227425
$data[3,14] = 0.028190957994264
$data[3,15] = 0.02819095799426401
$data[3,16] = 0.006953039255555555
$data[3,17] = 0.0625773533
$data[3,18] = 0.0001845033542644015
$data[3,19] = 0.0001845033542644015
# row 6: FAPs -> FAPs
$data[4,0] = 'FAPs'
$data[4,1] = 'Tnfsf10'
$data[4,2] = 'Tnfrsf11b'
$data[4,3] = 'FAPs'
$data[4,4] = 1
$data[4,5] = 0.3333333333333333
$data[4,6] = 0.09171866666666667
$data[4,7] = 0.275156
$data[4,8] = 0.006544770642485517
$data[4,9] = 0.006544770642485517
$data[4,10] = 3
$data[4,11] = 1
$data[4,12] = 2.095195666666667
$data[4,13] = 6.285587
$data[4,14] = 0.7791435378093522
$data[4,15] = 0.7791435378093522
$data[4,16] = 0.1921685529524444
$data[4,17] = 1.729516976572
$data[4,18] = 0.005099315752536953
$data[4,19] = 0.005099315752536953
# row 7: FAPs -> MuSCs
$data[5,0] = 'FAPs'
$data[5,1] = 'Tnfsf10'
$data[5,2] = 'Tnfrsf11b'
$data[5,3] = 'MuSCs'
$data[5,4] = 1
$data[5,5] = 0.3333333333333333
$data[5,6] = 0.09171866666666667
$data[5,7] = 0.275156
$data[5,8] = 0.006544770642485517
$data[5,9] = 0.006544770642485517
$data[5,10] = 3
$data[5,11] = 1
$data[5,12] = 0.5180969999999999
$data[5,13] = 1.554291
$data[5,14] = 0.1926655041963838
$data[5,15] = 0.1926655041963838
$data[5,16] = 0.04751916604399999
$data[5,17] = 0.427672494396
$data[5,18] = 0.001260951535684163
$data[5,19] = 0.001260951535684163
# row 8: MuSCs -> ECs
$data[6,0] = 'MuSCs'
$data[6,1] = 'Tnfsf10'
$data[6,2] = 'Tnfrsf11b'
$data[6,3] = 'ECs'
$data[6,4] = 1
$data[6,5] = 0.3333333333333333
$data[6,6] = 0.05557699999999999
$data[6,7] = 0.166731
$data[6,8] = 0.003965809046476372
$data[6,9] = 0.003965809046476372
$data[6,10] = 1
$data[6,11] = 0.3333333333333333
$data[6,12] = 0.07580833333333332
$data[6,13] = 0.227425
$data[6,14] = 0.028190957994264
$data[6,15] = 0.02819095799426401
$data[6,16] = 0.004213199741666666
$data[6,17] = 0.037918797675
$data[6,18] = 0.0001117999562424876
$data[6,19] = 0.0001117999562424876
# row 9: MuSCs -> FAPs
$data[7,0] = 'MuSCs'
$data[7,1] = 'Tnfsf10'
$data[7,2] = 'Tnfrsf11b'
$data[7,3] = 'FAPs'
$data[7,4] = 1
$data[7,5] = 0.3333333333333333
$data[7,6] = 0.05557699999999999
$data[7,7] = 0.166731
$data[7,8] = 0.003965809046476372
$data[7,9] = 0.003965809046476372
$data[7,10] = 3
$data[7,11] = 1
$data[7,12] = 2.095195666666667
$data[7,13] = 6.285587
$data[7,14] = 0.7791435378093522
$data[7,15] = 0.7791435378093522
$data[7,16] = 0.1164446895663333
$data[7,17] = 1.048002206097
$data[7,18] = 0.003089934490747934
$data[7,19] = 0.003089934490747934
# row 10: MuSCs -> MuSCs
$data[8,0] = 'MuSCs'
$data[8,1] = 'Tnfsf10'
$data[8,2] = 'Tnfrsf11b'
$data[8,3] = 'MuSCs'
$data[8,4] = 1
$data[8,5] = 0.3333333333333333
$data[8,6] = 0.05557699999999999
$data[8,7] = 0.166731
$data[8,8] = 0.003965809046476372
$data[8,9] = 0.003965809046476372
$data[8,10] = 3
$data[8,11] = 1
$data[8,12] = 0.5180969999999999
$data[8,13] = 1.554291
$data[8,14] = 0.1926655041963838
$data[8,15] = 0.1926655041963838
$data[8,16] = 0.02879427696899999
$data[8,17] = 0.259148492721
$data[8,18] = 0.0007640745994859504
$data[8,19] = 0.0007640745994859505

$ws.Range("A2:T10").Value = $data
